# Outstandings.xlsx maintenance edit:
#  - "Purchase 22-23": remove the trailing items 6/7/8 block (rows 31-36)
#    and drop the now-unused trailing "Chq no ..." column G.
#  - "Sale 22-23": the 09/23-24 ("b23-24MQ109") Putzmeister line is
#    replaced by the next one (b23-24MQ114, posted 45063, 369945) and the
#    intervening detail rows collapse into a single running-total row;
#    the running total formula now only covers the remaining rows.
#    Downstream "Sr. No" label shifts from 8 to 9.
#  - Selection / active-tab bookkeeping restored to "Purchase 22-23".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Purchase 22-23")
$ws2 = $wb.Worksheets.Item("Sale 22-23")

# --- Purchase 22-23 ------------------------------------------------------
# Drop rows 31-36 (blank separator + items 6, 7 and 8) entirely.
$ws1.Range("A31:A36").EntireRow.Delete()
# Drop the now-empty trailing "Chq no ..." column.
$ws1.Columns.Item(7).Delete()

# --- Sale 22-23 -----------------------------------------------------------
# Remove the old row 8 (b23-24MQ109) - row 9 (b23-24MQ114) slides up to
# become the new row 8.
$ws2.Range("A8:A8").EntireRow.Delete()
# Remove the remaining now-superseded detail rows (old rows 10-17).
$ws2.Range("A9:A16").EntireRow.Delete()
# Running total only needs to add up what is left (rows 5-8).
$ws2.Range("F8").Formula = "=E5+E6+E7+E8"
# The "Sr. No" label below shifts from 8 to 9.
$ws2.Range("A23").Value = 9

# --- View state -----------------------------------------------------------
[void]$ws2.Activate()
[void]$ws2.Range("A24").Select()

[void]$ws1.Activate()
[void]$ws1.Range("G16").Select()
